$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preprocess = 'convert to lower, trim "space" and ",", convert unicode to ascii, remove multiple spaces'
$features = '4 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone'
$model = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'
$templateFilter = '0 filters: '

$rows = @(
    @{ Row=32; Time='20160411_112803'; RunningTime=1038.199; TestAcc=0.973333333333333; ValAcc=0.993399339933993; J=0.175257731958763 },
    @{ Row=33; Time='20160411_114521'; RunningTime=1076.55;  TestAcc=0.980666666666667; ValAcc=0.986798679867987; J=0.157894736842105 },
    @{ Row=34; Time='20160411_120318'; RunningTime=1113.259; TestAcc=0.980666666666667; ValAcc=0.986798679867987; J=0.157894736842105 },
    @{ Row=35; Time='20160411_122151'; RunningTime=1159.053; TestAcc=0.982;              ValAcc=0.986798679867987; J=0.168421052631579 },
    @{ Row=36; Time='20160411_124110'; RunningTime=1162.679; TestAcc=0.979333333333333; ValAcc=0.993399339933993; J=0.154639175257732 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Time
    $ws.Cells.Item($row, 2).Value = $r.RunningTime
    $ws.Cells.Item($row, 3).Value = $preprocess
    $ws.Cells.Item($row, 4).Value = $features
    $ws.Cells.Item($row, 5).Value = $model
    $ws.Cells.Item($row, 6).Value = $modelDetails
    $ws.Cells.Item($row, 7).Value = $r.TestAcc
    $ws.Cells.Item($row, 8).Value = $r.ValAcc
    $ws.Cells.Item($row, 9).Value = $templateFilter
    $ws.Cells.Item($row, 10).Value = $r.J
}
